# Auto-generated edit script: updates market-board-derived profit columns (H:N)
# on the Sheets per the scheduled runner data refresh.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2525.5857
$ws.Range("I15").Value = 2525.5857
$ws.Range("K15").Value = 7576.757100000001
$ws.Range("M15").Value = -7407.757100000001
$ws.Range("H40").Value = 2335.3704
$ws.Range("I40").Value = 4433.6665
$ws.Range("J40").Value = 2073.0833
$ws.Range("K40").Value = 4433.6665
$ws.Range("L40").Value = 2073.0833
$ws.Range("M40").Value = -4258.6665
$ws.Range("N40").Value = -2423.0833
$ws.Range("H98").Value = 2007.5555
$ws.Range("I98").Value = 1672.1428
$ws.Range("J98").Value = 3181.5
$ws.Range("K98").Value = 1672.1428
$ws.Range("L98").Value = 3181.5
$ws.Range("M98").Value = -174.1428000000001
$ws.Range("N98").Value = -6177.5
$ws.Range("H121").Value = 2255.3333
$ws.Range("I121").Value = 3000
$ws.Range("J121").Value = 2187.6365
$ws.Range("K121").Value = 9000
$ws.Range("L121").Value = 6562.9095
$ws.Range("M121").Value = -7253
$ws.Range("N121").Value = -10056.9095
$ws.Range("H122").Value = 2007.5555
$ws.Range("I122").Value = 1672.1428
$ws.Range("J122").Value = 3181.5
$ws.Range("K122").Value = 5016.428400000001
$ws.Range("L122").Value = 9544.5
$ws.Range("M122").Value = -2566.428400000001
$ws.Range("N122").Value = -14444.5
$ws.Range("H132").Value = 7272.643
$ws.Range("I132").Value = 6649.2964
$ws.Range("J132").Value = 8394.666999999999
$ws.Range("K132").Value = 19947.8892
$ws.Range("L132").Value = 25184.001
$ws.Range("M132").Value = -17417.8892
$ws.Range("N132").Value = -30244.001
$ws.Range("H137").Value = 7464220
$ws.Range("I137").Value = 11365201
$ws.Range("J137").Value = 1473.1305
$ws.Range("K137").Value = 34095603
$ws.Range("L137").Value = 4419.3915
$ws.Range("M137").Value = -34093053
$ws.Range("N137").Value = -9519.3915
$ws.Range("H138").Value = 1992.591
$ws.Range("I138").Value = 2160.25
$ws.Range("J138").Value = 1852.875
$ws.Range("K138").Value = 6480.75
$ws.Range("L138").Value = 5558.625
$ws.Range("M138").Value = -1340.75
$ws.Range("N138").Value = -15838.625
$ws.Range("H141").Value = 6018.12
$ws.Range("I141").Value = 2145.5652
$ws.Range("J141").Value = 50552.5
$ws.Range("K141").Value = 6436.6956
$ws.Range("L141").Value = 151657.5
$ws.Range("M141").Value = -1256.6956
$ws.Range("N141").Value = -162017.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5140763.5
$ws.Range("I32").Value = 6904.877
$ws.Range("K32").Value = 6904.877
$ws.Range("M32").Value = -6617.877
$ws.Range("H132").Value = 1896997.9
$ws.Range("I132").Value = 1343.2094
$ws.Range("J132").Value = 5441047.5
$ws.Range("K132").Value = 4029.6282
$ws.Range("L132").Value = 16323142.5
$ws.Range("M132").Value = -1499.6282
$ws.Range("N132").Value = -16328202.5

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 747
$ws.Range("I94").Value = 747
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 747
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -296
$ws.Range("N94").ClearContents()
$ws.Range("H134").Value = 3429.3044
$ws.Range("I134").Value = 1245.3478
$ws.Range("J134").Value = 7797.2173
$ws.Range("K134").Value = 3736.0434
$ws.Range("L134").Value = 23391.6519
$ws.Range("M134").Value = -1201.0434
$ws.Range("N134").Value = -28461.6519

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4832762
$ws.Range("J31").Value = 9806095
$ws.Range("L31").Value = 9806095
$ws.Range("N31").Value = -9806685
$ws.Range("H34").Value = 4832762
$ws.Range("J34").Value = 9806095
$ws.Range("L34").Value = 9806095
$ws.Range("N34").Value = -9806499
$ws.Range("H134").Value = 2718.4
$ws.Range("I134").Value = 1256.8572
$ws.Range("J134").Value = 3997.25
$ws.Range("K134").Value = 3770.5716
$ws.Range("L134").Value = 11991.75
$ws.Range("M134").Value = -1235.5716
$ws.Range("N134").Value = -17061.75

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 5250
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 5250
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 15750
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -16098
$ws.Range("H129").Value = 2620.2354
$ws.Range("I129").Value = 571.8
$ws.Range("J129").Value = 3473.75
$ws.Range("K129").Value = 1715.4
$ws.Range("L129").Value = 10421.25
$ws.Range("M129").Value = 3284.6
$ws.Range("N129").Value = -20421.25

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2979246.2
$ws.Range("I102").Value = 6495979.5
$ws.Range("J102").Value = 3549.2307
$ws.Range("K102").Value = 6495979.5
$ws.Range("L102").Value = 3549.2307
$ws.Range("M102").Value = -6494357.5
$ws.Range("N102").Value = -6793.2307
$ws.Range("H122").Value = 1252938.1
$ws.Range("I122").Value = 2501404.2
$ws.Range("J122").Value = 4472
$ws.Range("K122").Value = 7504212.600000001
$ws.Range("L122").Value = 13416
$ws.Range("M122").Value = -7501762.600000001
$ws.Range("N122").Value = -18316

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 8546.091
$ws.Range("I13").Value = 2250
$ws.Range("J13").Value = 25335.666
$ws.Range("K13").Value = 2250
$ws.Range("L13").Value = 25335.666
$ws.Range("M13").Value = -2110
$ws.Range("N13").Value = -25615.666
$ws.Range("H68").Value = 3299.5334
$ws.Range("I68").Value = 3137.1428
$ws.Range("J68").Value = 3441.625
$ws.Range("K68").Value = 3137.1428
$ws.Range("L68").Value = 3441.625
$ws.Range("M68").Value = -2388.1428
$ws.Range("N68").Value = -4939.625
$ws.Range("H71").Value = 3299.5334
$ws.Range("I71").Value = 3137.1428
$ws.Range("J71").Value = 3441.625
$ws.Range("K71").Value = 15685.714
$ws.Range("L71").Value = 17208.125
$ws.Range("M71").Value = -11941.714
$ws.Range("N71").Value = -24696.125

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 881.4483
$ws.Range("I113").Value = 470.42856
$ws.Range("J113").Value = 1012.2273
$ws.Range("K113").Value = 1411.28568
$ws.Range("L113").Value = 3036.6819
$ws.Range("M113").Value = 758.71432
$ws.Range("N113").Value = -7376.6819
$ws.Range("H126").Value = 7439.68
$ws.Range("I126").Value = 8371.048000000001
$ws.Range("J126").Value = 2550
$ws.Range("K126").Value = 25113.144
$ws.Range("L126").Value = 7650
$ws.Range("M126").Value = -22643.144
$ws.Range("N126").Value = -12590
